$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.641.39'
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").Value = '2.960.44'
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '496.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.37%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -5.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.15'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.49%  '
$ws.Range("E10").Value = '  -6.47%  '
$ws.Range("E11").Value = '  -5.32%  '
$ws.Range("D12").Value = '3.465.97'
$ws.Range("E13").Value = '  -3.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.50%  '
$ws.Range("E15").Value = '  -8.94%  '
$ws.Range("D16").Value = '56.667.09'
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.964.83'
$ws.Range("E17").Value = '  -3.81%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '317.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.23%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("E24").Value = '  -3.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E27").Value = '  -4.70%  '
$ws.Range("D28").Value = '0.0₃0863'
$ws.Range("E28").Value = '  -12.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.88%  '
$ws.Range("E30").Value = '  -5.94%  '
$ws.Range("E31").Value = '  -6.31%  '
$ws.Range("E32").Value = '  -6.16%  '
$ws.Range("E33").Value = '  -8.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.30%  '
$ws.Range("E36").Value = '  -5.26%  '
$ws.Range("E37").Value = '  -9.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.98%  '
$ws.Range("E39").Value = '  -7.30%  '
$ws.Range("D40").Value = '2.990.89'
$ws.Range("E40").Value = '  -3.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("E43").Value = '  -4.14%  '
$ws.Range("E44").Value = '  -7.36%  '
$ws.Range("D45").Value = '2.153.21'
$ws.Range("E45").Value = '  -8.31%  '
$ws.Range("E46").Value = '  -9.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.85'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.923'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -11.92%  '
$ws.Range("E49").Value = '  -5.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.06%  '
$ws.Range("E51").Value = '  -13.97%  '
